# Applies the July 10, 2020 data refresh to the COVID disparities sheet.
# Updates Date Published, Total Cases/Deaths, Black/AA case & death counts,
# and the derived percentages for the affected state/county rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Texas
$ws.Range("B2").Value = 44022
$ws.Range("C2").Value = 24758
$ws.Range("D2").Value = 710
$ws.Range("E2").Value = 2851
$ws.Range("F2").Value = 91
$ws.Range("G2").Value = 11.52
$ws.Range("H2").Value = 12.82

# Row 3 - Texas -- Bexar County
$ws.Range("B3").Value = 44022
$ws.Range("C3").Value = 18602
$ws.Range("D3").Value = 166

# Row 8 - Kentucky
$ws.Range("B8").Value = 44022
$ws.Range("C8").Value = 18670
$ws.Range("D8").Value = 620
$ws.Range("E8").Value = 1756
$ws.Range("F8").Value = 88
$ws.Range("G8").Value = 13.57
$ws.Range("H8").Value = 15.29
$ws.Range("K8").Value = 12943
$ws.Range("L8").Value = 581

# Row 9 - Arkansas (Date Published unchanged, already 44022)
$ws.Range("C9").Value = 26803
$ws.Range("D9").Value = 313
$ws.Range("E9").Value = 5796
$ws.Range("G9").Value = 25.01
$ws.Range("H9").Value = 25.26
$ws.Range("K9").Value = 23171
$ws.Range("L9").Value = 293

# Row 10 - California - San Diego
$ws.Range("B10").Value = 44022
$ws.Range("C10").Value = 18863
$ws.Range("D10").Value = 420
$ws.Range("E10").Value = 690
$ws.Range("G10").Value = 4.58
$ws.Range("H10").Value = 3.89
$ws.Range("K10").Value = 15058
$ws.Range("L10").Value = 411

# Row 12 - New Mexico
$ws.Range("B12").Value = 44022
$ws.Range("C12").Value = 14549
$ws.Range("D12").Value = 539
$ws.Range("E12").Value = 277
$ws.Range("G12").Value = 1.9

# Row 16 - California - Los Angeles
$ws.Range("B16").Value = 44021
$ws.Range("C16").Value = 127358
$ws.Range("D16").Value = 3738
$ws.Range("E16").Value = 3407
$ws.Range("F16").Value = 380
$ws.Range("G16").Value = 4.76
$ws.Range("H16").Value = 10.91
$ws.Range("K16").Value = 71628
$ws.Range("L16").Value = 3482

# Row 24 - Colorado
$ws.Range("B24").Value = 44022
$ws.Range("C24").Value = 36191
$ws.Range("D24").Value = 1724
$ws.Range("E24").Value = 1871
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 6.35
$ws.Range("H24").Value = 7.04
$ws.Range("K24").Value = 29452
$ws.Range("L24").Value = 1661

# Row 25 - Nebraska
$ws.Range("B25").Value = 44022
$ws.Range("C25").Value = 20777
$ws.Range("D25").Value = 286
$ws.Range("E25").Value = 1230
$ws.Range("G25").Value = 7.61
$ws.Range("H25").Value = 8.03
$ws.Range("K25").Value = 16169
$ws.Range("L25").Value = 274

# Row 32 - Washington
$ws.Range("B32").Value = 44022
$ws.Range("C32").Value = 39218
$ws.Range("D32").Value = 1424
$ws.Range("E32").Value = 1520
$ws.Range("F32").Value = 46
$ws.Range("G32").Value = 5.43
$ws.Range("H32").Value = 3.39
$ws.Range("K32").Value = 27979
$ws.Range("L32").Value = 1355

# Row 36 - Iowa (Date Published, Total Deaths, Count Deaths Black/AA,
# and Pct Deaths Black/AA unchanged)
$ws.Range("C36").Value = 34172
$ws.Range("E36").Value = 2983
$ws.Range("G36").Value = 8.73

# Row 39 - Idaho
$ws.Range("B39").Value = 44022
$ws.Range("C39").Value = 9928
$ws.Range("D39").Value = 101
$ws.Range("E39").Value = 145
$ws.Range("G39").Value = 1.46
$ws.Range("H39").Value = 0.99
